$d = $word.ActiveDocument

$tab = [char]9

# ---------------------------------------------------------------------
# 1) The "{{ clients[0].signature }} ... Date: {{ signature_date }}" line
#    becomes "{{ clients[0].name }} ... Date: {{ signature_date }}"
#    (the underlined "Date: " also loses its underline formatting).
# ---------------------------------------------------------------------
$p7 = $d.Paragraphs(7)
$rng7 = $p7.Range
$old1 = "{{ clients[0].signature }}              " + $tab + $tab
$new1 = "{{ clients[0].name }}" + $tab + $tab + $tab
$rng7.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2)

# Clear the leftover underline formatting that used to decorate "Date: "
$p7 = $d.Paragraphs(7)
$rng7b = $p7.Range
$rng7b.Find.ClearFormatting()
$found = $rng7b.Find.Execute("Date: ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rng7b.Font.Underline = 0
}

# ---------------------------------------------------------------------
# 2) The "Signature" label line becomes the signature placeholder.
# ---------------------------------------------------------------------
$p8 = $d.Paragraphs(8)
$rng8 = $p8.Range
$old2 = "Signature" + $tab + "      "
$new2 = "{{ clients[0].signature }}              "
$rng8.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $new2, 2)

# ---------------------------------------------------------------------
# 3) First-page footer: " " + "tdd: 617.371.1228" merge into one run.
# ---------------------------------------------------------------------
$sec = $d.Sections(1)
$ftr = $sec.Footers(2)
$rngF = $ftr.Range
$oldF = " tdd: 617.371.1228"
$rngF.Find.Execute($oldF, $true, $false, $false, $false, $false, $true, 1, $false, $oldF, 2)
